# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and swap the WrappedBTC/WrappedEther rows (16-17) per upstream ranking change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.240.71'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '3.404.58'
$ws.Range("E3").Value = '  -3.50%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'583.24"
$ws.Range("E5").Value = '  -2.62%  '
$ws.Range("D6").Value = "'136.45"
$ws.Range("E6").Value = '  -5.15%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.400.11'
$ws.Range("E8").Value = '  -3.67%  '
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").Value = "'7.20"
$ws.Range("E10").Value = '  -7.43%  '
$ws.Range("E12").Value = '  -7.83%  '
$ws.Range("D13").Value = '3.986.08'
$ws.Range("E13").Value = '  -3.57%  '
$ws.Range("D14").Value = "'0.0000177"
$ws.Range("E14").Value = '  -11.21%  '
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.409.59'
$ws.Range("E16").Value = '  -3.00%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '65.246.01'
$ws.Range("E17").Value = '  -1.00%  '
$ws.Range("D18").Value = "'25.87"
$ws.Range("E18").Value = '  -9.74%  '
$ws.Range("D19").Value = "'9.71"
$ws.Range("E19").Value = '  -11.10%  '
$ws.Range("D20").Value = "'5.85"
$ws.Range("E20").Value = '  -5.52%  '
$ws.Range("D21").Value = "'13.50"
$ws.Range("E21").Value = '  -5.77%  '
$ws.Range("D22").Value = "'383.20"
$ws.Range("E23").Value = '  -7.36%  '
$ws.Range("D24").Value = "'72.57"
$ws.Range("E24").Value = '  -6.23%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '3.544.27'
$ws.Range("E26").Value = '  -3.52%  '
$ws.Range("D27").Value = "'0.0000105"
$ws.Range("E27").Value = '  -10.01%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").Value = "'7.04"
$ws.Range("E29").Value = '  -9.47%  '
$ws.Range("E30").Value = '  -9.99%  '
$ws.Range("D31").Value = "'8.03"
$ws.Range("D32").Value = '3.411.99'
$ws.Range("E32").Value = '  -3.33%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -8.07%  '
$ws.Range("D35").Value = "'22.69"
$ws.Range("E35").Value = '  -6.84%  '
$ws.Range("D36").Value = "'170.63"
$ws.Range("E36").Value = '  -2.67%  '
$ws.Range("D37").Value = "'6.73"
$ws.Range("E37").Value = '  -10.42%  '
$ws.Range("D38").Value = "'1.44"
$ws.Range("E38").Value = '  -8.47%  '
$ws.Range("E39").Value = '  -12.84%  '
$ws.Range("E40").Value = '  -10.51%  '
$ws.Range("D41").Value = "'0.0754"
$ws.Range("E41").Value = '  -7.76%  '
$ws.Range("D42").Value = "'0.814"
$ws.Range("E42").Value = '  -5.12%  '
$ws.Range("D43").Value = "'43.46"
$ws.Range("E43").Value = '  -4.05%  '
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("E45").Value = '  -13.64%  '
$ws.Range("E46").Value = '  -11.19%  '
$ws.Range("D47").Value = "'1.09"
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").Value = "'22.07"
$ws.Range("E48").Value = '  -2.16%  '
$ws.Range("E49").Value = '  -8.31%  '
$ws.Range("D50").Value = "'2.02"
$ws.Range("E50").Value = '  -15.83%  '
$ws.Range("D51").Value = '2.171.42'
$ws.Range("E51").Value = '  -7.60%  '
